$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New backlog item: unfinished user story from Sprint 1 retrospective
$ws.Range("A19").Value = "Jako zarządca budynku mogę otrzymać informację o pomieszczeniach w budynku, które przekraczają określony poziom zużycia energii cieplnej / m^3 podany jako parametr, aby znaleźć miejsca do poprawy w infrastrukturze."
$ws.Range("C19").Value = 4

$ws.Range("B20").Value = "Zaimplementowanie poszukiwania pomieszczeń przekraczającyh określony poziom zużycia energi cieplnej"
$ws.Range("C20").Value = 3

$ws.Range("B21").Value = "Pobieranie podanej przez administratora wartości energii"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "0.5"
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C21").PasteSpecial(-4122) | Out-Null

$ws.Range("B22").Value = "Wyświetlenie znalezionych pomieszczeń"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0.5"
$ws.Range("C20").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

# New backlog item: composite pattern fix
$ws.Range("A23").Value = "Poprawienie wzorca kopmpozyt"
$ws.Range("C23").Value = 1

$ws.Range("B24").Value = "Poprawienie funkcji tak aby operowały na abstrakcji"
$ws.Range("C24").Value = 1

# Row heights to match content wrapping (Excel would normally auto-fit these)
$ws.Rows.Item(19).RowHeight = 78.75
$ws.Rows.Item(20).RowHeight = 31.5

# View state
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F26").Select()
